$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row of data
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "GH-BLUE"
$ws.Range("C7").Value = "G Handbag"
$ws.Range("D7").Value = 200

# Update selection to E7
$ws.Range("E7").Select()

# Resize the workbook window
$excel.Width = 26840
$excel.Height = 13320
